$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row (row 1) with three more labelled columns.
$ws.Range("B1").Value = "Hi"
$ws.Range("C1").Value = "a"
$ws.Range("D1").Value = "b"

# Column E: filled down with "c", with occasional "cc" entries,
# matching the pattern found in the edited workbook. Row 21 is written
# and then cleared below, so it ends up with no data at all.
$cValues = @(
  "c",  # 1
  "c",  # 2
  "c",  # 3
  "c",  # 4
  "c",  # 5
  "c",  # 6
  "c",  # 7
  "cc", # 8
  "c",  # 9
  "c",  # 10
  "c",  # 11
  "c",  # 12
  "cc", # 13
  "c",  # 14
  "c",  # 15
  "c",  # 16
  "c",  # 17
  "c",  # 18
  "c",  # 19
  "cc", # 20
  "cc", # 21
  "cc", # 22
  "cc", # 23
  "cc"  # 24
)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 5).Value = $cValues[$i]
}

# Row 21 ends up completely blank (no cell data at all) in the final
# workbook, so clear what was just written there.
$ws.Range("E21").ClearContents() | Out-Null

$ws.Range("G7").Select() | Out-Null
